# "fixed export and fixing maps"
# The workbook held an export with an extra subtitle row and two extra
# census-year columns (1989 / 2002). This reverts it to the simpler,
# single-year (2014) export and restores the sheet's real name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from the placeholder "1" to the municipality name.
$ws.Name = "ჩოხატაური"

# Drop the subtitle row ("(მოსახლეობის აღწერის შედეგებით)") - row 2 -
# shifting everything below it up by one.
$ws.Rows(2).Delete()

# Drop the 1989 and 2002 columns (old C:D), keeping only the 2014 column
# that is now column B.
$ws.Columns("C:D").Delete()

# The remaining year/value column now holds the 2014 figures.
$ws.Range("B4").Value = 2014
$ws.Range("B5").Value = 825.1

# Clear the now-empty formatted-but-blank cells left over from the old
# layout so they drop out of the sheet entirely.
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()

# Leave the selection where the author left it when they saved.
$ws.Range("A2").Select() | Out-Null
